$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.836.76"
$ws.Range("E2").Value = "  -3.23%  "

$ws.Range("D3").Value = "2.275.32"
$ws.Range("E3").Value = "  -3.80%  "

$ws.Range("E4").Value = "  +0.10%  "

$r = $ws.Range("D5")
$r.Value = "'531.03"
$r.ClearFormats()
$ws.Range("E5").Value = "  -4.69%  "

$r = $ws.Range("D6")
$r.Value = "'130.64"
$r.ClearFormats()
$ws.Range("E6").Value = "  -1.86%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("E8").Value = "  +0.48%  "

$ws.Range("D9").Value = "2.275.10"
$ws.Range("E9").Value = "  -3.67%  "

$r = $ws.Range("D10")
$r.Value = "'0.0989"
$r.ClearFormats()
$ws.Range("E10").Value = "  -5.59%  "

$r = $ws.Range("D11")
$r.Value = "'5.43"
$r.ClearFormats()
$ws.Range("E11").Value = "  -3.41%  "

$ws.Range("E12").Value = "  -0.29%  "

$r = $ws.Range("D13")
$r.Value = "'0.329"
$r.ClearFormats()
$ws.Range("E13").Value = "  -3.66%  "

$r = $ws.Range("D14")
$r.Value = "'23.42"
$r.ClearFormats()
$ws.Range("E14").Value = "  -3.23%  "

$ws.Range("D15").Value = "2.686.24"
$ws.Range("E15").Value = "  -3.47%  "

$ws.Range("D16").Value = "57.843.36"
$ws.Range("E16").Value = "  -3.15%  "

$ws.Range("E17").Value = "  -4.27%  "

$ws.Range("D18").Value = "2.276.25"
$ws.Range("E18").Value = "  -3.93%  "

$r = $ws.Range("D19")
$r.Value = "'10.47"
$r.ClearFormats()
$ws.Range("E19").Value = "  -5.25%  "

$r = $ws.Range("D20")
$r.Value = "'4.21"
$r.ClearFormats()
$ws.Range("E20").Value = "  -5.51%  "

$r = $ws.Range("D21")
$r.Value = "'311.36"
$r.ClearFormats()
$ws.Range("E21").Value = "  -2.63%  "

$r = $ws.Range("D22")
$r.Value = "'6.36"
$r.ClearFormats()
$ws.Range("E22").Value = "  -4.11%  "

$ws.Range("E23").Value = "  +0.00%  "

$r = $ws.Range("D24")
$r.Value = "'62.56"
$r.ClearFormats()
$ws.Range("E24").Value = "  -2.52%  "

$ws.Range("E25").Value = "  -2.66%  "

$r = $ws.Range("D26")
$r.Value = "'1.00"
$r.ClearFormats()
$ws.Range("E26").Value = "  +0.03%  "

$r = $ws.Range("D27")
$r.Value = "'7.96"
$r.ClearFormats()
$ws.Range("E27").Value = "  -4.84%  "

$r = $ws.Range("D28")
$r.Value = "'1.27"
$r.ClearFormats()
$ws.Range("E28").Value = "  -6.93%  "

$r = $ws.Range("D29")
$r.Value = "'170.30"
$r.ClearFormats()
$ws.Range("E29").Value = "  -0.24%  "

$r = $ws.Range("D30")
$r.Value = "'1.69"
$r.ClearFormats()
$ws.Range("E30").Value = "  -6.04%  "

$ws.Range("D31").Value = "0.0₃0716"
$ws.Range("E31").Value = "  -5.22%  "

$r = $ws.Range("D32")
$r.Value = "'5.73"
$r.ClearFormats()
$ws.Range("E32").Value = "  -5.49%  "

$r = $ws.Range("D33")
$r.Value = "'1.04"
$r.ClearFormats()
$ws.Range("E33").Value = "  -5.61%  "

$r = $ws.Range("D34")
$r.Value = "'0.378"
$r.ClearFormats()
$ws.Range("E34").Value = "  -4.81%  "

$r = $ws.Range("D36")
$r.Value = "'17.69"
$r.ClearFormats()
$ws.Range("E36").Value = "  -2.19%  "

$ws.Range("E37").Value = "  -0.04%  "

$r = $ws.Range("D38")
$r.Value = "'1.22"
$r.ClearFormats()
$ws.Range("E38").Value = "  -6.74%  "

$r = $ws.Range("D39")
$r.Value = "'3.88"
$r.ClearFormats()
$ws.Range("E39").Value = "  -5.57%  "

$r = $ws.Range("D40")
$r.Value = "'38.16"
$r.ClearFormats()
$ws.Range("E40").Value = "  -1.11%  "

$r = $ws.Range("D41")
$r.Value = "'1.48"
$r.ClearFormats()
$ws.Range("E41").Value = "  -6.30%  "

$r = $ws.Range("D42")
$r.Value = "'141.33"
$r.ClearFormats()
$ws.Range("E42").Value = "  -2.10%  "

$r = $ws.Range("D43")
$r.Value = "'286.39"
$r.ClearFormats()
$ws.Range("E43").Value = "  -9.67%  "

$r = $ws.Range("D44")
$r.Value = "'3.40"
$r.ClearFormats()
$ws.Range("E44").Value = "  -3.35%  "

$r = $ws.Range("D45")
$r.Value = "'0.0948"
$r.ClearFormats()
$ws.Range("E45").Value = "  -1.60%  "

$r = $ws.Range("D46")
$r.Value = "'0.0493"
$r.ClearFormats()
$ws.Range("E46").Value = "  -3.13%  "

$r = $ws.Range("D47")
$r.Value = "'0.548"
$r.ClearFormats()
$ws.Range("E47").Value = "  -3.29%  "

$r = $ws.Range("D48")
$r.Value = "'18.02"
$r.ClearFormats()
$ws.Range("E48").Value = "  -6.70%  "

$r = $ws.Range("D49")
$r.Value = "'0.0209"
$r.ClearFormats()
$ws.Range("E49").Value = "  -3.58%  "

$ws.Range("E50").Value = "  -1.14%  "
